$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from H1 (bold/bordered header style) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-51
$iValues = @(6, 5, 6, 5, 7, 8, 7, 7, 6, 7, 7, 8, 7, 7, 7, 7, 6, 7, 7, 7, 6, 6, 7, 7, 7, 8, 8, 6, 8, 6, 7, 7, 7, 9, 5, 6, 4, 6, 6, 8, 7, 6, 6, 6, 5, 4, 4, 4, 5, 4)
$jValues = @(7, 6, 6, 6, 7, 8, 8, 7, 6, 8, 8, 8, 7, 8, 7, 7, 7, 8, 8, 7, 7, 7, 8, 8, 8, 8, 8, 7, 8, 7, 8, 7, 8, 9, 6, 8, 6, 6, 6, 8, 7, 7, 8, 7, 5, 5, 6, 6, 6, 4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}